$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 115, shifting existing rows 115:223 down to 116:224
$ws.Rows.Item(115).Insert()

$ws.Cells.Item(115, 1).Value = 10
$ws.Cells.Item(115, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(115, 3).Value = "La Araucanía"
$ws.Cells.Item(115, 4).Value = 44484
$ws.Cells.Item(115, 5).Value = 9
$ws.Cells.Item(115, 6).Value = "Fruta"
$ws.Cells.Item(115, 7).Value = 100108
$ws.Cells.Item(115, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(115, 9).Value = 100108002
$ws.Cells.Item(115, 10).Value = "Mango"
$ws.Cells.Item(115, 11).Value = "Sin especificar"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 285
$ws.Cells.Item(115, 14).Value = 8000
$ws.Cells.Item(115, 15).Value = 9000
$ws.Cells.Item(115, 16).Value = 8386
$ws.Cells.Item(115, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(115, 18).Value = "Perú"
$ws.Cells.Item(115, 19).Value = 2096
$ws.Cells.Item(115, 20).Value = 4
